$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# "3. " + "Einbeziehung" + " von Personen" (3 runs) -> "3. Einbeziehung von Personen" (1 run)
# Word's Find/Replace merges the touched run with its identically-formatted
# neighbours, so replacing just the middle run's text with itself triggers
# the merge into a single run while leaving the separate leading "  " run
# (the list-number indent spacer) untouched.
$d.Content.Find.Execute("Einbeziehung", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Einbeziehung", 2) | Out-Null

# --- Change 2 -----------------------------------------------------------
# "2. Spezifisch, Messbar, " + "Attraktiv" + ", Realistisch, T" + "ypisch"
# (4 runs) -> "2. Spezifisch, Messbar, Attraktiv, Realistisch, Typisch" (1 run)
$d.Content.Find.Execute("Attraktiv", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Attraktiv", 2) | Out-Null

# --- Change 3 -----------------------------------------------------------
# Fix the marked-correct answer for question 14 (the SMART acronym
# question) from "2" to "1": the run " - 2" must become two runs,
# " - " and "1", each keeping the original <w:lang w:val="de-DE"/> rPr.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Frage 14: Anwort(en) - 2`r") {
        $e = $p.Range.End

        # Insert the replacement digit "1" immediately before the old "2".
        $insertionPoint = $d.Range($e - 2, $e - 2)
        $insertionPoint.InsertBefore("1")

        # The old "2" now sits one character later; delete it.
        $old2 = $d.Range($e - 1, $e)
        $old2.Delete()

        # The new "1" is currently merged into the preceding " - " run
        # (same formatting). Toggle Bold on/off on just that character to
        # force Word to split it into its own run while keeping the
        # inherited rPr (lang de-DE) explicit on the new run.
        $newDigit = $d.Range($e - 2, $e - 1)
        $newDigit.Bold = $true
        $newDigit.Bold = $false

        break
    }
}
